$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D152").Value = 44476
$ws.Range("D153").Value = 44476
$ws.Range("D154").Value = 44386
$ws.Range("D155").Value = 44386
$ws.Range("D156").Value = 44306
$ws.Range("D157").Value = 44306
$ws.Range("D158").Value = 44425
$ws.Range("D159").Value = 44425
$ws.Range("D160").Value = 44187
$ws.Range("D161").Value = 44187

# New row 162 - Primera
$ws.Cells.Item(162, 1).Value = 11
$ws.Cells.Item(162, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(162, 3).Value = "Bíobío"
$ws.Cells.Item(162, 4).Value = 44250
$ws.Cells.Item(162, 5).Value = 8
$ws.Cells.Item(162, 6).Value = 100112009
$ws.Cells.Item(162, 7).Value = "Acelga"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 200
$ws.Cells.Item(162, 11).Value = 600
$ws.Cells.Item(162, 12).Value = 700
$ws.Cells.Item(162, 13).Value = 650
$ws.Cells.Item(162, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(162, 15).Value = "Región de Ñuble"
$ws.Cells.Item(162, 16).Value = 650
$ws.Cells.Item(162, 17).Value = 1
$ws.Cells.Item(162, 18).Value = "Hortaliza"

# New row 163 - Segunda
$ws.Cells.Item(163, 1).Value = 11
$ws.Cells.Item(163, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(163, 3).Value = "Bíobío"
$ws.Cells.Item(163, 4).Value = 44250
$ws.Cells.Item(163, 5).Value = 8
$ws.Cells.Item(163, 6).Value = 100112009
$ws.Cells.Item(163, 7).Value = "Acelga"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Segunda"
$ws.Cells.Item(163, 10).Value = 100
$ws.Cells.Item(163, 11).Value = 500
$ws.Cells.Item(163, 12).Value = 500
$ws.Cells.Item(163, 13).Value = 500
$ws.Cells.Item(163, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(163, 15).Value = "Región de Ñuble"
$ws.Cells.Item(163, 16).Value = 500
$ws.Cells.Item(163, 17).Value = 1
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# Apply date formatting to D162:D163 to match D column formatting
$ws.Range("D162").NumberFormat = $ws.Range("D161").NumberFormat
$ws.Range("D163").NumberFormat = $ws.Range("D161").NumberFormat
